$d = $word.ActiveDocument

function New-WordXmlPackage($innerBodyXml) {
    $header = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>'
    $footer = '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    return $header + $innerBodyXml + $footer
}

# --- Paragraph 1: "Playlist (create, findByUsr, findById, Add)" ---
#   -> "Playlist (create, findByUsr, findById)"  (remove the grey ", Add" run,
#      split findByUsr/findById into their own runs bracketed by proofErr marks)
$playlistPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Playlist (create,")) {
        $playlistPara = $p
        break
    }
}
if ($null -eq $playlistPara) {
    throw "Could not locate the 'Playlist (create, ...)' paragraph"
}

$full = $d.Range($playlistPara.Range.Start, $playlistPara.Range.End - 1)

$inner = ''
$inner += '<w:r w:rsidRPr="00C96025"><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">Playlist (create, </w:t></w:r>'
$inner += '<w:proofErr w:type="spellStart"/>'
$inner += '<w:r w:rsidRPr="00C96025"><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>findByUsr</w:t></w:r>'
$inner += '<w:proofErr w:type="spellEnd"/>'
$inner += '<w:r w:rsidRPr="00C96025"><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r>'
$inner += '<w:proofErr w:type="spellStart"/>'
$inner += '<w:r w:rsidRPr="00C96025"><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>findById</w:t></w:r>'
$inner += '<w:proofErr w:type="spellEnd"/>'
$inner += '<w:r w:rsidRPr="00C96025"><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="en-GB"/></w:rPr><w:t>)</w:t></w:r>'

$full.InsertXML((New-WordXmlPackage $inner))

# --- Paragraph 2: "Track (create, findById)" ---
#   -> "Track (create, findById, Add)"
$trackPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Track (create,")) {
        $trackPara = $p
        break
    }
}
if ($null -eq $trackPara) {
    throw "Could not locate the 'Track (create, ...)' paragraph"
}

$full2 = $d.Range($trackPara.Range.Start, $trackPara.Range.End - 1)

$inner2 = ''
$inner2 += '<w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">Track (create, </w:t></w:r>'
$inner2 += '<w:proofErr w:type="spellStart"/>'
$inner2 += '<w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="en-GB"/></w:rPr><w:t>findById</w:t></w:r>'
$inner2 += '<w:proofErr w:type="spellEnd"/>'
$inner2 += '<w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="en-GB"/></w:rPr><w:t>, Add</w:t></w:r>'
$inner2 += '<w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="en-GB"/></w:rPr><w:t>)</w:t></w:r>'

$full2.InsertXML((New-WordXmlPackage $inner2))

Write-Host "Edit complete"
